$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference row used as the source of the existing table border/format
# (the full thin box border currently used by the data rows A2:E9).
$srcFormat = $ws.Range("A9")

# Build the desired "left + right thin border only" look once, on an
# out-of-the-way scratch cell, so it can be stamped onto both A10 and
# C10 via a simple format copy (keeps a single shared cell style for
# both, same as applying it once and filling across).
$scratch = $ws.Range("Z1")
$srcFormat.Copy()
$scratch.PasteSpecial(-4122)             # xlPasteFormats
$scratch.Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
$scratch.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none

# --- New row 10: MathNet.Numerics entry -------------------------------

# Column A (Toolkit)
$a10 = $ws.Range("A10")
$a10.Value = "MathNet.Numerics"
$scratch.Copy()
$a10.PasteSpecial(-4122)          # xlPasteFormats

# Column C (Licence)
$c10 = $ws.Range("C10")
$c10.Value = "MIT/X11"
$scratch.Copy()
$c10.PasteSpecial(-4122)          # xlPasteFormats

# Column D (Link) - plain, unstyled cell
$d10 = $ws.Range("D10")
$d10.Value = "https://numerics.mathdotnet.com/License.html"

# Clean up the scratch cell so it doesn't leave stray content/format behind.
$scratch.Clear()
$excel.CutCopyMode = 0

# --- View state ---------------------------------------------------------
$null = $ws.Range("E16").Select()
